# Update dashboards - 2025-10-31
# Refresh the "Latest Date" and Present/Lag1-4 series for the rows whose
# underlying FRED data rolled forward: T5YIFR (29), T10YIE (30), DGS2 (48),
# DGS5 (49), DGS10 (50), MORTGAGE30US (51) and DBAA (52).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateText {
    param($range, [string]$text)
    # Leading apostrophe forces the engine to keep the value as literal
    # text instead of auto-converting a date-shaped string into a date
    # serial number (matches how the source file stores these as
    # explicit text, not formatted dates).
    $ws.Range($range).Value = "'" + $text
}

# Row 29 - T5YIFR
Set-DateText "N29" "2025-10-30"
$ws.Range("Q29").Value = 2.2
$ws.Range("R29").Value = 2.3
$ws.Range("S29").Value = 2.35
$ws.Range("T29").Value = 2.32
$ws.Range("U29").Value = 2.27

# Row 30 - T10YIE
Set-DateText "N30" "2025-10-30"
$ws.Range("Q30").Value = 2.29
$ws.Range("R30").Value = 2.3
$ws.Range("S30").Value = 2.28
$ws.Range("T30").Value = 2.28
$ws.Range("U30").Value = 2.29

# Row 48 - DGS2
Set-DateText "N48" "2025-10-29"
$ws.Range("Q48").Value = 3.59
$ws.Range("R48").Value = 3.47
$ws.Range("S48").Value = 3.48
$ws.Range("T48").Value = 3.48
$ws.Range("U48").Value = 3.48

# Row 49 - DGS5
Set-DateText "N49" "2025-10-29"
$ws.Range("Q49").Value = 3.7
$ws.Range("R49").Value = 3.6
$ws.Range("S49").Value = 3.61
$ws.Range("T49").Value = 3.61
$ws.Range("U49").Value = 3.61

# Row 50 - DGS10
Set-DateText "N50" "2025-10-29"
$ws.Range("Q50").Value = 4.08
$ws.Range("R50").Value = 3.99
$ws.Range("S50").Value = 4.01
$ws.Range("T50").Value = 4.02
$ws.Range("U50").Value = 4.01

# Row 51 - MORTGAGE30US
Set-DateText "N51" "2025-10-30"
$ws.Range("Q51").Value = 6.17
$ws.Range("R51").Value = 6.19
$ws.Range("S51").Value = 6.27
$ws.Range("T51").Value = 6.3
$ws.Range("U51").Value = 6.34

# Row 52 - DBAA
Set-DateText "N52" "2025-10-29"
$ws.Range("Q52").Value = 5.69
$ws.Range("R52").Value = 5.64
$ws.Range("S52").Value = 5.64
$ws.Range("T52").Value = 5.67
$ws.Range("U52").Value = 5.67
